$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 1575
$ws.Range("D11").Value = 1575
$ws.Range("D12").Value = 1421.172084527004
$ws.Range("D13").Value = 1421.172084527004
